$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the current row 273, shifting existing rows
# (273-325) down to (277-329). The inserted rows inherit the row-273
# formatting (column D keeps its date style).
$ws.Range("A273:A276").EntireRow.Insert()

# Common values shared by the 4 new "Murcott" records (date 2021-09-24 /
# serial 44463) being inserted ahead of the existing data block.
$rows = @(273, 274, 275, 276)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = 8
    $ws.Cells.Item($r, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44463
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = "Cítricos"
    $ws.Cells.Item($r, 9).Value = 100102004
    $ws.Cells.Item($r, 10).Value = "Mandarina"
    $ws.Cells.Item($r, 11).Value = "Murcott"
    $ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"
    $ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 20).Value = 10
}

# Row 273: Especial
$ws.Cells.Item(273, 12).Value = "Especial"
$ws.Cells.Item(273, 13).Value = 400
$ws.Cells.Item(273, 14).Value = 5000
$ws.Cells.Item(273, 15).Value = 5500
$ws.Cells.Item(273, 16).Value = 5250
$ws.Cells.Item(273, 19).Value = 525

# Row 274: Primera
$ws.Cells.Item(274, 12).Value = "Primera"
$ws.Cells.Item(274, 13).Value = 400
$ws.Cells.Item(274, 14).Value = 4000
$ws.Cells.Item(274, 15).Value = 4500
$ws.Cells.Item(274, 16).Value = 4250
$ws.Cells.Item(274, 19).Value = 425

# Row 275: Segunda
$ws.Cells.Item(275, 12).Value = "Segunda"
$ws.Cells.Item(275, 13).Value = 360
$ws.Cells.Item(275, 14).Value = 3000
$ws.Cells.Item(275, 15).Value = 3500
$ws.Cells.Item(275, 16).Value = 3250
$ws.Cells.Item(275, 19).Value = 325

# Row 276: Tercera
$ws.Cells.Item(276, 12).Value = "Tercera"
$ws.Cells.Item(276, 13).Value = 200
$ws.Cells.Item(276, 14).Value = 2000
$ws.Cells.Item(276, 15).Value = 2500
$ws.Cells.Item(276, 16).Value = 2250
$ws.Cells.Item(276, 19).Value = 225
